$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H62").Value = 4466.8335
$ws.Range("I62").Value = 4466.8335
$ws.Range("K62").Value = 4466.8335
$ws.Range("M62").Value = -3842.8335
$ws.Range("H65").Value = 4466.8335
$ws.Range("I65").Value = 4466.8335
$ws.Range("K65").Value = 22334.1675
$ws.Range("M65").Value = -19214.1675
$ws.Range("H69").Value = 25000
$ws.Range("J69").Value = 25000
$ws.Range("L69").Value = 75000
$ws.Range("N69").Value = -76748
$ws.Range("H72").Value = 25000
$ws.Range("J72").Value = 25000
$ws.Range("L72").Value = 225000
$ws.Range("N72").Value = -233736
$ws.Range("H86").Value = 70179510
$ws.Range("H89").Value = 70179510
$ws.Range("H98").Value = 3835.8438
$ws.Range("I98").Value = 2850.3684
$ws.Range("J98").Value = 5276.154
$ws.Range("K98").Value = 2850.3684
$ws.Range("L98").Value = 5276.154
$ws.Range("M98").Value = -1352.3684
$ws.Range("N98").Value = -8272.154
$ws.Range("H122").Value = 3835.8438
$ws.Range("I122").Value = 2850.3684
$ws.Range("J122").Value = 5276.154
$ws.Range("K122").Value = 8551.1052
$ws.Range("L122").Value = 15828.462
$ws.Range("M122").Value = -6101.1052
$ws.Range("N122").Value = -20728.462

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 870.4583
$ws.Range("I2").Value = 870.4583
$ws.Range("J2").Value = 0
$ws.Range("K2").Value = 870.4583
$ws.Range("L2").Value = 0
$ws.Range("M2").Value = -757.4583
$ws.Range("N2").ClearContents()
$ws.Range("H32").Value = 252780.55
$ws.Range("I32").Value = 295462.72
$ws.Range("J32").Value = 10915
$ws.Range("K32").Value = 295462.72
$ws.Range("L32").Value = 10915
$ws.Range("M32").Value = -295175.72
$ws.Range("N32").Value = -11489
$ws.Range("H45").Value = 46614
$ws.Range("I45").Value = 61608.47
$ws.Range("J45").Value = 4129.6665
$ws.Range("K45").Value = 61608.47
$ws.Range("L45").Value = 4129.6665
$ws.Range("M45").Value = -61231.47
$ws.Range("N45").Value = -4883.6665
$ws.Range("H63").Value = 23397.684
$ws.Range("I63").Value = 7581
$ws.Range("J63").Value = 30697.691
$ws.Range("K63").Value = 7581
$ws.Range("L63").Value = 30697.691
$ws.Range("M63").Value = -6895
$ws.Range("N63").Value = -32069.691
$ws.Range("H66").Value = 23397.684
$ws.Range("I66").Value = 7581
$ws.Range("J66").Value = 30697.691
$ws.Range("K66").Value = 37905
$ws.Range("L66").Value = 153488.455
$ws.Range("M66").Value = -34473
$ws.Range("N66").Value = -160352.455
$ws.Range("H116").Value = 870.4583
$ws.Range("I116").Value = 870.4583
$ws.Range("J116").Value = 0
$ws.Range("K116").Value = 870.4583
$ws.Range("L116").Value = 0
$ws.Range("M116").Value = 1423.5417
$ws.Range("N116").ClearContents()
$ws.Range("H132").Value = 2320.0833
$ws.Range("I132").Value = 2346.356
$ws.Range("J132").Value = 770
$ws.Range("K132").Value = 7039.068000000001
$ws.Range("L132").Value = 2310
$ws.Range("M132").Value = -4509.068000000001
$ws.Range("N132").Value = -7370

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 870.4583
$ws.Range("I3").Value = 870.4583
$ws.Range("J3").Value = 0
$ws.Range("K3").Value = 870.4583
$ws.Range("L3").Value = 0
$ws.Range("M3").Value = -756.4583
$ws.Range("N3").ClearContents()
$ws.Range("H105").Value = 11719.417
$ws.Range("I105").Value = 15948.143
$ws.Range("K105").Value = 15948.143
$ws.Range("M105").Value = -14201.143
$ws.Range("H134").Value = 19566668
$ws.Range("I134").Value = 1437.8572
$ws.Range("K134").Value = 4313.571599999999
$ws.Range("M134").Value = -1778.571599999999

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H10").Value = 4776.6
$ws.Range("I10").Value = 4007
$ws.Range("K10").Value = 4007
$ws.Range("M10").Value = -3868
$ws.Range("H31").Value = 3156.442
$ws.Range("I31").Value = 2881.6956
$ws.Range("K31").Value = 2881.6956
$ws.Range("M31").Value = -2586.6956
$ws.Range("H34").Value = 3156.442
$ws.Range("I34").Value = 2881.6956
$ws.Range("K34").Value = 2881.6956
$ws.Range("M34").Value = -2679.6956
$ws.Range("H134").Value = 1336.3846
$ws.Range("I134").Value = 1193.3889
$ws.Range("K134").Value = 3580.1667
$ws.Range("M134").Value = -1045.1667

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 1100391.2
$ws.Range("I4").Value = 1292005.9
$ws.Range("K4").Value = 3876017.7
$ws.Range("M4").Value = -3875905.7
$ws.Range("H97").Value = 322.16666
$ws.Range("I97").Value = 410.75
$ws.Range("K97").Value = 1232.25
$ws.Range("M97").Value = -736.25
$ws.Range("H137").Value = 4054.2727
$ws.Range("I137").Value = 4363
$ws.Range("J137").Value = 3797
$ws.Range("K137").Value = 13089
$ws.Range("L137").Value = 11391
$ws.Range("M137").Value = -7989
$ws.Range("N137").Value = -21591

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H3").Value = 3199.8
$ws.Range("J3").Value = 2749.75
$ws.Range("L3").Value = 2749.75
$ws.Range("N3").Value = -2981.75
$ws.Range("H10").Value = 7399.4
$ws.Range("J10").Value = 9331.333000000001
$ws.Range("L10").Value = 9331.333000000001
$ws.Range("N10").Value = -9669.333000000001
$ws.Range("H11").Value = 1679999.9
$ws.Range("I11").Value = 1679999.9
$ws.Range("K11").Value = 1679999.9
$ws.Range("M11").Value = -1679860.9
$ws.Range("H12").Value = 5799.5
$ws.Range("I12").Value = 1600
$ws.Range("K12").Value = 1600
$ws.Range("M12").Value = -1460
$ws.Range("H97").Value = 701.4167
$ws.Range("I97").Value = 706.8
$ws.Range("J97").Value = 674.5
$ws.Range("K97").Value = 706.8
$ws.Range("L97").Value = 674.5
$ws.Range("M97").Value = -210.8
$ws.Range("N97").Value = -1666.5
$ws.Range("H107").Value = 77347
$ws.Range("J107").Value = 499.125
$ws.Range("L107").Value = 499.125
$ws.Range("N107").Value = -4339.125
$ws.Range("H132").Value = 574112.2
$ws.Range("I132").Value = 5739.577
$ws.Range("K132").Value = 17218.731
$ws.Range("M132").Value = -14688.731

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 2109.4285
$ws.Range("I40").Value = 1312.4546
$ws.Range("K40").Value = 1312.4546
$ws.Range("M40").Value = -1176.4546
$ws.Range("H54").Value = 28333.334
$ws.Range("I54").Value = 0
$ws.Range("J54").Value = 28333.334
$ws.Range("K54").Value = 0
$ws.Range("L54").Value = 28333.334
$ws.Range("M54").ClearContents()
$ws.Range("N54").Value = -29621.334
$ws.Range("H61").Value = 2515.9644
$ws.Range("I61").Value = 2363.3845
$ws.Range("K61").Value = 2363.3845
$ws.Range("M61").Value = -2161.3845
$ws.Range("H74").Value = 46749.25
$ws.Range("I74").Value = 46999.5
$ws.Range("J74").Value = 46665.832
$ws.Range("K74").Value = 46999.5
$ws.Range("L74").Value = 46665.832
$ws.Range("M74").Value = -46001.5
$ws.Range("N74").Value = -48661.832
$ws.Range("H77").Value = 46749.25
$ws.Range("I77").Value = 46999.5
$ws.Range("J77").Value = 46665.832
$ws.Range("K77").Value = 140998.5
$ws.Range("L77").Value = 139997.496
$ws.Range("M77").Value = -136006.5
$ws.Range("N77").Value = -149981.496
$ws.Range("H100").Value = 2841.4546
$ws.Range("I100").Value = 2875.1
$ws.Range("J100").Value = 2505
$ws.Range("K100").Value = 2875.1
$ws.Range("L100").Value = 2505
$ws.Range("M100").Value = -2334.1
$ws.Range("N100").Value = -3587
$ws.Range("H113").Value = 2515.9644
$ws.Range("I113").Value = 2363.3845
$ws.Range("K113").Value = 2363.3845
$ws.Range("M113").Value = -193.3845000000001
$ws.Range("H122").Value = 3088.3726
$ws.Range("I122").Value = 2601.0286
$ws.Range("J122").Value = 4154.4375
$ws.Range("K122").Value = 7803.085800000001
$ws.Range("L122").Value = 12463.3125
$ws.Range("M122").Value = -5353.085800000001
$ws.Range("N122").Value = -17363.3125
$ws.Range("H132").Value = 3125.9062
$ws.Range("I132").Value = 2847.2727
$ws.Range("J132").Value = 3738.9
$ws.Range("K132").Value = 8541.8181
$ws.Range("L132").Value = 11216.7
$ws.Range("M132").Value = -6011.8181
$ws.Range("N132").Value = -16276.7

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H21").Value = 0
$ws.Range("I21").Value = 0
$ws.Range("K21").Value = 0
$ws.Range("M21").ClearContents()
$ws.Range("H35").Value = 0
$ws.Range("I35").Value = 0
$ws.Range("K35").Value = 0
$ws.Range("M35").ClearContents()
$ws.Range("H81").Value = 51629.2
$ws.Range("I81").Value = 1887.1875
$ws.Range("K81").Value = 3774.375
$ws.Range("M81").Value = -2713.375
$ws.Range("H84").Value = 51629.2
$ws.Range("I84").Value = 1887.1875
$ws.Range("K84").Value = 18871.875
$ws.Range("M84").Value = -13567.875
$ws.Range("H132").Value = 2110.423
$ws.Range("I132").Value = 1758.75
$ws.Range("J132").Value = 2673.1
$ws.Range("K132").Value = 5276.25
$ws.Range("L132").Value = 8019.299999999999
$ws.Range("M132").Value = -2746.25
$ws.Range("N132").Value = -13079.3
$ws.Range("H136").Value = 34419.355
$ws.Range("I136").Value = 48418.715
$ws.Range("K136").Value = 145256.145
$ws.Range("M136").Value = -142706.145
